# feat: add 2022-Q1 data
#
# 1) Insert a new worksheet "2022-Q1" right before the "总计" (total) sheet,
#    holding the per-fund holdings snapshot for 2022-Q1 (same layout as the
#    other quarterly sheets).
# 2) Prepend a "2022-Q1" summary row to the "总计" sheet (date / count /
#    total market value), shifting the existing rows down and renumbering
#    the leading index column.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper: write a text-typed value into a cell without leaving a residual
# non-default style behind (mirrors the "@" / text cells used for the
# percentage-ish figures in the source data, e.g. "32.13").
# ---------------------------------------------------------------------
function Set-TextCell($cell, [string]$text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

# ---------------------------------------------------------------------
# 1. Create the "2022-Q1" sheet, positioned right before "总计"
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$q1 = $wb.Worksheets.Add($totalSheet)
$q1.Name = "2022-Q1"

# Template sheet to borrow header / index-column formatting from.
$template = $wb.Worksheets.Item("2021-Q4")

# Header row (B1:H1) formatting + index column (A2:A6) formatting.
$template.Range("B1:H1").Copy() | Out-Null
$q1.Range("B1:H1").PasteSpecial(-4122) | Out-Null

$template.Range("A2:A6").Copy() | Out-Null
$q1.Range("A2:A6").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false

$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

$q1.Range("A2").Value = 0
Set-TextCell $q1.Range("B2") "001838"
$q1.Range("C2").Value = "国投瑞银国家安全灵活配置混合"
Set-TextCell $q1.Range("D2") "32.13"
Set-TextCell $q1.Range("E2") "94.68"
Set-TextCell $q1.Range("F2") "7.20"
Set-TextCell $q1.Range("G2") "2.3134"
$q1.Range("H2").Value = 7

$q1.Range("A3").Value = 1
Set-TextCell $q1.Range("B3") "005774"
$q1.Range("C3").Value = "华夏产业升级混合"
Set-TextCell $q1.Range("D3") "12.64"
Set-TextCell $q1.Range("E3") "94.47"
Set-TextCell $q1.Range("F3") "8.29"
Set-TextCell $q1.Range("G3") "1.0479"
$q1.Range("H3").Value = 4

$q1.Range("A4").Value = 2
Set-TextCell $q1.Range("B4") "460002"
$q1.Range("C4").Value = "华泰柏瑞积极成长混合A"
Set-TextCell $q1.Range("D4") "6.11"
Set-TextCell $q1.Range("E4") "81.55"
Set-TextCell $q1.Range("F4") "8.91"
Set-TextCell $q1.Range("G4") "0.5444"
$q1.Range("H4").Value = 3

$q1.Range("A5").Value = 3
Set-TextCell $q1.Range("B5") "009317"
$q1.Range("C5").Value = "金信核心竞争力灵活配置混合"
Set-TextCell $q1.Range("D5") "0.19"
Set-TextCell $q1.Range("E5") "89.48"
Set-TextCell $q1.Range("F5") "2.88"
Set-TextCell $q1.Range("G5") "0.0055"
$q1.Range("H5").Value = 9

$q1.Range("A6").Value = 4
Set-TextCell $q1.Range("B6") "960030"
$q1.Range("C6").Value = "华泰柏瑞积极成长混合H"
Set-TextCell $q1.Range("D6") "0.00"
Set-TextCell $q1.Range("E6") "81.55"
Set-TextCell $q1.Range("F6") "8.91"
$q1.Range("G6").Value = 0
$q1.Range("H6").Value = 3

# ---------------------------------------------------------------------
# 2. Prepend the 2022-Q1 summary row to "总计", pushing the rest down and
#    renumbering the A-column running index (0,1,2,...).
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

$rows = @(
    @("2022-Q1", 5,  3.91),
    @("2021-Q4", 6,  5.75),
    @("2021-Q3", 9,  10.92),
    @("2021-Q2", 12, 9.08),
    @("2021-Q1", 13, 9.4),
    @("2020-Q4", 47, 24.44)
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $total.Range("A$r").Value = $i
    $total.Range("B$r").Value = $rows[$i][0]
    $total.Range("C$r").Value = $rows[$i][1]
    $total.Range("D$r").Value = $rows[$i][2]
}

# Row 7 is brand new (the sheet previously only had 6 rows), so its A-cell
# has no inherited "index column" formatting yet - copy it over from A6.
$total.Range("A6").Copy() | Out-Null
$total.Range("A7").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Restore the original active sheet/selection (adding a sheet makes it the
# active one) so the view state is left as it was before the edit.
$wb.Worksheets.Item("2020-Q4").Activate()
